# Update Price (col D) and Volume(1h) (col E) figures in the cryptos list.
# Numeric-looking Price values are prefixed with a leading apostrophe so
# Excel keeps them as text (matching the original inline-string cells)
# instead of silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.314.18'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '2.445.91'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''584.30'
$ws.Range('E5').Value = '  +2.20%  '
$ws.Range('D6').Value = '''143.74'
$ws.Range('E6').Value = '  -2.09%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('D9').Value = '2.443.93'
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('E10').Value = '  -3.13%  '
$ws.Range('D11').Value = '''0.160'
$ws.Range('E11').Value = '  +2.08%  '
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').Value = '''0.345'
$ws.Range('E13').Value = '  -3.14%  '
$ws.Range('D14').Value = '''26.46'
$ws.Range('E14').Value = '  -2.15%  '
$ws.Range('E15').Value = '  -3.66%  '
$ws.Range('D16').Value = '2.874.81'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').Value = '62.067.17'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('D18').Value = '2.425.94'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('E19').Value = '  -3.70%  '
$ws.Range('E20').Value = '  -2.79%  '
$ws.Range('D21').Value = '''330.56'
$ws.Range('E21').Value = '  +0.82%  '
$ws.Range('E22').Value = '  -1.91%  '
$ws.Range('E23').Value = '  -6.00%  '
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('D25').Value = '''65.87'
$ws.Range('E25').Value = '  +0.83%  '
$ws.Range('D26').Value = '''9.39'
$ws.Range('E26').Value = '  +4.31%  '
$ws.Range('D27').Value = '''619.94'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Value = '2.568.45'
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('D29').Value = '0.0₃0957'
$ws.Range('E29').Value = '  -7.39%  '
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('E31').Value = '  -4.70%  '
$ws.Range('E32').Value = '  -3.47%  '
$ws.Range('E33').Value = '  +0.26%  '
$ws.Range('E34').Value = '  -0.79%  '
$ws.Range('D35').Value = '''4.93'
$ws.Range('E35').Value = '  -5.79%  '
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('E37').Value = '  -6.14%  '
$ws.Range('D38').Value = '''0.376'
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('D39').Value = '''151.71'
$ws.Range('E39').Value = '  +4.25%  '
$ws.Range('D40').Value = '''18.35'
$ws.Range('E40').Value = '  -2.18%  '
$ws.Range('D41').Value = '''5.25'
$ws.Range('E41').Value = '  -3.36%  '
$ws.Range('D42').Value = '''1.76'
$ws.Range('E42').Value = '  -1.72%  '
$ws.Range('D43').Value = '''42.46'
$ws.Range('E43').Value = '  +1.32%  '
$ws.Range('D45').Value = '''2.47'
$ws.Range('E45').Value = '  -8.67%  '
$ws.Range('D46').Value = '''143.48'
$ws.Range('E46').Value = '  -3.66%  '
$ws.Range('E47').Value = '  -3.39%  '
$ws.Range('D48').Value = '''0.0526'
$ws.Range('D49').Value = '''0.599'
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('D50').Value = '''19.51'
$ws.Range('E50').Value = '  -7.98%  '
$ws.Range('E51').Value = '  -1.23%  '
